# ---------------------------------------------------------------------------
# feat: add 2022-Q1 data
#   - insert a new "2022-Q1" worksheet (fund holdings detail) right before
#     the "总计" (totals) summary sheet
#   - prepend a new "2022-Q1" row to the "总计" summary table
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item("总计")

# 1. Create the new sheet immediately before "总计" and name it.
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# 2. Header row, styled like the other quarterly sheets (bold/border/center).
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $newSheet.Cells.Item(1, $col).Value = $headers[$col - 2]
}

# 3. Fund holdings detail rows.
#    Columns D (基金规模), E (股票总仓位), F (仓位占比), G (持有市值) are kept
#    as literal text (matching the source data, which stores them as
#    inline strings, not numbers - e.g. so values like "21.89" are not
#    silently renormalised and codes that happen to look numeric keep any
#    leading zeros). Column H (仓位排名) is a genuine number.
$rows = @(
    @("001822", "华商智能生活灵活配置混合", "21.89", "90.29", "6.78", "1.4841", 4),
    @("166301", "华商新趋势优选灵活配置混合", "26.96", "86.39", "2.60", "0.7010", 6),
    @("630002", "华商盛世成长混合", "22.62", "93.39", "2.82", "0.6379", 9),
    @("001933", "华商新兴活力灵活配置混合", "6.91", "92.28", "7.03", "0.4858", 4),
    @("000390", "华商优势行业混合", "19.61", "88.06", "2.39", "0.4687", 9),
    @("001959", "华商乐享互联灵活配置混合", "3.74", "88.84", "3.23", "0.1208", 6),
    @("004044", "金鹰转型动力灵活配置混合", "0.72", "93.34", "6.15", "0.0443", 5),
    @("008488", "华商恒益稳健混合", "2.03", "58.93", "1.61", "0.0327", 9),
    @("010756", "兴华永兴混合A", "0.35", "94.57", "4.54", "0.0159", 5),
    @("540007", "汇丰晋信中小盘股票", "0.61", "93.28", "2.14", "0.0131", 6),
    @("002303", "金鹰智慧生活灵活配置混合", "0.11", "89.88", "6.91", "0.0076", 3),
    @("010999", "兴华瑞丰混合A", "0.06", "29.21", "3.52", "0.0021", 4),
    @("011000", "兴华瑞丰混合C", "0.05", "29.21", "3.52", "0.0018", 4),
    @("010757", "兴华永兴混合C", "0.01", "94.57", "4.54", "0.0005", 5)
)

$lastRow = 1 + $rows.Count

# Pre-format B:G of the data area as Text so numeric-looking values (fund
# codes with leading zeros, percentages, NAV figures) are stored as strings
# rather than being auto-coerced into numbers - this reuses a single shared
# style for the whole block instead of minting one per cell.
$newSheet.Range("B2:G" + $lastRow).NumberFormat = "@"

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = $r - 2
    $newSheet.Cells.Item($r, 2).Value = $row[0]
    $newSheet.Cells.Item($r, 3).Value = $row[1]
    $newSheet.Cells.Item($r, 4).Value = $row[2]
    $newSheet.Cells.Item($r, 5).Value = $row[3]
    $newSheet.Cells.Item($r, 6).Value = $row[4]
    $newSheet.Cells.Item($r, 7).Value = $row[5]
    $newSheet.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# 4. Apply the bold/border/center "index" style (already used for column A
#    and the header row elsewhere in the workbook) to the header row and to
#    column A of the data rows, by copying formats from an existing cell
#    that already carries it - this guarantees we land on the exact same
#    shared style instead of minting a near-duplicate.
$totalSheet.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$totalSheet.Range("A2").Copy()
$newSheet.Range("A2:A" + $lastRow).PasteSpecial(-4122)

$excel.CutCopyMode = $false
